# Update automàtic: dades i banners [2026-02-25 18:50]
# Applies refreshed meteocat data values to the summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-25 18:48:27'
$ws.Range("K2").Value = '12.9 MJ/m2'
$ws.Range("E3").Value = '2026-02-25 18:48:29'
$ws.Range("E4").Value = '2026-02-25 18:48:31'
$ws.Range("J4").Value = '1021.6 hPa'
$ws.Range("O4").Value = '8.3 °C'
$ws.Range("E5").Value = '2026-02-25 18:48:34'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '26%'
$ws.Range("N5").Value = '3.0 °C 18:16 TU'
$ws.Range("E6").Value = '2026-02-25 18:48:37'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '89%'
$ws.Range("E7").Value = '2026-02-25 18:48:39'
$ws.Range("E8").Value = '2026-02-25 18:48:42'
$ws.Range("O8").Value = '12.1 °C'
$ws.Range("E9").Value = '2026-02-25 18:48:44'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '90%'
$ws.Range("O9").Value = '10.0 °C'
$ws.Range("E10").Value = '2026-02-25 18:48:46'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '88%'
$ws.Range("O10").Value = '9.7 °C'
$ws.Range("E11").Value = '2026-02-25 18:48:47'
$ws.Range("E12").Value = '2026-02-25 18:48:48'
$ws.Range("O12").Value = '9.9 °C'
$ws.Range("E13").Value = '2026-02-25 18:48:49'
$ws.Range("K13").Value = '15.1 MJ/m2'
$ws.Range("O13").Value = '6.6 °C'
$ws.Range("E14").Value = '2026-02-25 18:48:50'
$ws.Range("E15").Value = '2026-02-25 18:48:52'
$ws.Range("E16").Value = '2026-02-25 18:48:55'
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = '29%'
$ws.Range("N16").Value = '1.6 °C 18:27 TU'
$ws.Range("E17").Value = '2026-02-25 18:48:58'
$ws.Range("N17").Value = '6.0 °C 18:29 TU'
$ws.Range("O17").Value = '9.6 °C'
$ws.Range("E18").Value = '2026-02-25 18:49:00'
$ws.Range("J18").Value = '1021.7 hPa'
$ws.Range("E19").Value = '2026-02-25 18:49:03'
$ws.Range("E20").Value = '2026-02-25 18:49:06'
$ws.Range("N20").Value = '0.0 °C 18:22 TU'
$ws.Range("O20").Value = '3.0 °C'
$ws.Range("E21").Value = '2026-02-25 18:49:08'
$ws.Range("J21").Value = '1021.3 hPa'
$ws.Range("O21").Value = '10.0 °C'
$ws.Range("E22").Value = '2026-02-25 18:49:11'
$ws.Range("E23").Value = '2026-02-25 18:49:14'
$ws.Range("E24").Value = '2026-02-25 18:49:16'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '76%'
$ws.Range("J24").Value = '1019.8 hPa'
$ws.Range("E25").Value = '2026-02-25 18:49:19'
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = '32%'
$ws.Range("O25").Value = '5.5 °C'
$ws.Range("E26").Value = '2026-02-25 18:49:21'
$ws.Range("J26").Value = '1019.0 hPa'
$ws.Range("N26").Value = '7.1 °C 18:23 TU'
$ws.Range("O26").Value = '10.6 °C'
$ws.Range("E27").Value = '2026-02-25 18:49:24'
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = '40%'
$ws.Range("E28").Value = '2026-02-25 18:49:27'
$ws.Range("E29").Value = '2026-02-25 18:49:29'
$ws.Range("E30").Value = '2026-02-25 18:49:31'
$ws.Range("E31").Value = '2026-02-25 18:49:34'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '92%'
$ws.Range("E32").Value = '2026-02-25 18:49:37'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '50%'
$ws.Range("E33").Value = '2026-02-25 18:49:39'
$ws.Range("J33").Value = '1020.9 hPa'
$ws.Range("O33").Value = '8.6 °C'
$ws.Range("E34").Value = '2026-02-25 18:49:42'
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = '50%'
$ws.Range("O34").Value = '3.8 °C'
$ws.Range("E35").Value = '2026-02-25 18:49:45'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '36%'
$ws.Range("J35").Value = '1019.1 hPa'
$ws.Range("O35").Value = '12.8 °C'
$ws.Range("E36").Value = '2026-02-25 18:49:47'
$ws.Range("E37").Value = '2026-02-25 18:49:50'
$ws.Range("E38").Value = '2026-02-25 18:49:53'
$ws.Range("O38").Value = '9.2 °C'
$ws.Range("E39").Value = '2026-02-25 18:49:55'
$ws.Range("O39").Value = '2.6 °C'
$ws.Range("E40").Value = '2026-02-25 18:49:58'
$ws.Range("E41").Value = '2026-02-25 18:50:00'
$ws.Range("E42").Value = '2026-02-25 18:50:03'
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = '92%'
$ws.Range("O42").Value = '11.5 °C'
$ws.Range("E43").Value = '2026-02-25 18:50:05'
$ws.Range("E44").Value = '2026-02-25 18:50:08'
$ws.Range("E45").Value = '2026-02-25 18:50:10'
$ws.Range("J45").Value = '1019.5 hPa'
$ws.Range("K45").Value = '13.8 MJ/m2'
$ws.Range("L45").Value = '20.9 km/h - 119º 18:28 TU'
$ws.Range("E46").Value = '2026-02-25 18:50:13'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '81%'
$ws.Range("J46").Value = '1020.5 hPa'
$ws.Range("O46").Value = '9.5 °C'
